$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "29.071.52"
$ws.Range("E2").Value = "  -0.04%  "

# Row 3
Set-TextCell "D3" "1.836.62"
$ws.Range("E3").Value = "  +0.48%  "

# Row 4
Set-TextCell "D4" "1.000"
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
Set-TextCell "D5" "243.05"
$ws.Range("E5").Value = "  +0.63%  "

# Row 6
Set-TextCell "D6" "0.6289"
$ws.Range("E6").Value = "  -1.12%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
Set-TextCell "D8" "0.07569"
$ws.Range("E8").Value = "  +2.98%  "

# Row 9
Set-TextCell "D9" "0.2938"
$ws.Range("E9").Value = "  +0.04%  "

# Row 10
Set-TextCell "D10" "22.63"
$ws.Range("E10").Value = "  -0.85%  "

# Row 11
Set-TextCell "D11" "0.07750"
$ws.Range("E11").Value = "  +0.98%  "

# Row 12
Set-TextCell "D12" "1.834.08"
$ws.Range("E12").Value = "  +0.33%  "

# Row 13
Set-TextCell "D13" "4.971"
$ws.Range("E13").Value = "  -0.43%  "

# Row 14
Set-TextCell "D14" "0.6665"
$ws.Range("E14").Value = "  +0.31%  "

# Row 15
Set-TextCell "D15" "0.00001001"
$ws.Range("E15").Value = "  +15.23%  "

# Row 16
Set-TextCell "D16" "83.14"
$ws.Range("E16").Value = "  +1.25%  "

# Row 17
Set-TextCell "D17" "6.084"

# Row 18
Set-TextCell "D18" "29.096.84"
$ws.Range("E18").Value = "  +0.08%  "

# Row 19
Set-TextCell "D19" "226.70"
$ws.Range("E19").Value = "  -0.18%  "

# Row 20
$ws.Range("E20").Value = "  +0.32%  "

# Row 21
Set-TextCell "D21" "1.001"
$ws.Range("E21").Value = "  +0.12%  "

# Row 22
Set-TextCell "D22" "7.230"
$ws.Range("E22").Value = "  +1.32%  "

# Row 23
$ws.Range("E23").Value = "  +0.07%  "

# Row 24
Set-TextCell "D24" "159.86"
$ws.Range("E24").Value = "  +0.91%  "

# Row 25
$ws.Range("E25").Value = "  +1.68%  "

# Row 26
Set-TextCell "D26" "8.506"
$ws.Range("E26").Value = "  +0.24%  "

# Row 27
Set-TextCell "D27" "17.93"
$ws.Range("E27").Value = "  +0.04%  "

# Row 28
Set-TextCell "D28" "1.497"
$ws.Range("E28").Value = "  -0.44%  "

# Row 29
$ws.Range("E29").Value = "  +0.18%  "

# Row 30
Set-TextCell "D30" "4.012"
$ws.Range("E30").Value = "  -0.50%  "

# Row 31
Set-TextCell "D31" "1.194"
$ws.Range("E31").Value = "  -0.78%  "

# Row 32
Set-TextCell "D32" "0.05257"
$ws.Range("E32").Value = "  -1.35%  "

# Row 33
Set-TextCell "D33" "1.852"
$ws.Range("E33").Value = "  +0.64%  "

# Row 34
Set-TextCell "D34" "0.7380"
$ws.Range("E34").Value = "  +0.03%  "

# Row 35
$ws.Range("E35").Value = "  -1.69%  "

# Row 36
Set-TextCell "D36" "2.680"
$ws.Range("E36").Value = "  +1.15%  "

# Row 37
Set-TextCell "D37" "1.245.74"
$ws.Range("E37").Value = "  -4.32%  "

# Row 38
Set-TextCell "D38" "2.763"
$ws.Range("E38").Value = "  +0.69%  "

# Row 39
Set-TextCell "D39" "0.01786"
$ws.Range("E39").Value = "  -0.11%  "

# Row 40
Set-TextCell "D40" "6.389"
$ws.Range("E40").Value = "  +1.19%  "

# Row 41
Set-TextCell "D41" "0.9015"
$ws.Range("E41").Value = "  -0.12%  "

# Row 42
Set-TextCell "D42" "1.001"
$ws.Range("E42").Value = "  +0.19%  "

# Row 43
$ws.Range("B43").Value = "BabyDogeCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D43" "0.00000000130"
$ws.Range("E43").Value = "  +8.97%  "

# Row 44
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D44" "102.13"
$ws.Range("E44").Value = "  -0.58%  "

# Row 45
Set-TextCell "D45" "1.990.30"
$ws.Range("E45").Value = "  +0.77%  "

# Row 46
Set-TextCell "D46" "64.41"
$ws.Range("E46").Value = "  +0.45%  "

# Row 47
$ws.Range("E47").Value = "  -0.21%  "

# Row 48
$ws.Range("E48").Value = "  +1.44%  "

# Row 49
Set-TextCell "D49" "8.944"
$ws.Range("E49").Value = "  +2.15%  "

# Row 50
Set-TextCell "D50" "0.05773"
$ws.Range("E50").Value = "  -0.73%  "

# Row 51
Set-TextCell "D51" "6.723"
$ws.Range("E51").Value = "  +0.32%  "
